# correção das notas do fórum para matc65 em 2021.2
# For every data row (2..50), if any of the daily-view columns (B..H)
# contains a 1 (a recorded forum view), zero out the whole row's
# view-tracking columns (B..H) as well as the derived totals
# (I = total_views, J = nota_view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $hasView = $false
    for ($col = 2; $col -le 8; $col++) {  # B..H
        $val = $ws.Cells.Item($row, $col).Value2
        if ($val -eq 1) {
            $hasView = $true
        }
    }

    if ($hasView) {
        for ($col = 2; $col -le 10; $col++) {  # B..J
            $ws.Cells.Item($row, $col).Value = 0
        }
    }
}
